# "Create Survey updated by team 3"
# Adds three new field types (Satisfactory, Image Radio, Image Checkbox) to
# the "For every case check those items:" guideline block and to the
# Type -> type-string mapping table at the bottom of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Make room: insert 11 new rows before the existing "Type" header
#    (old row 23), pushing it (and everything below) down to row 34.
# ---------------------------------------------------------------------
$ws.Rows("23:33").Insert()

# ---------------------------------------------------------------------
# Helper style constants re-used throughout (match existing sheet look):
#   column A label cells -> bold, 12pt, Calibri   (style "s=4")
#   column C value cells -> regular, 12pt, Calibri (style "s=3")
# ---------------------------------------------------------------------

function Set-Label($addr, $text) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 12
    $c.Font.Bold = $true
    $c.Value = $text
    $ws.Range($addr).EntireRow.RowHeight = 15.75
}

function Set-Value($addr, $text) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 12
    $c.Font.Bold = $false
    $c.Value = $text
    $ws.Range($addr).EntireRow.RowHeight = 15.75
}

function Set-NoteRedStar($addr, $starText, $restText) {
    # Two-colour note: a red "***" (or "*** ") marker followed by the
    # regular-coloured explanation text, as plain runs inside one cell.
    $c = $ws.Range($addr)
    $c.Font.Name = "Calibri"
    $c.Font.Size = 12
    $c.Value = $starText + $restText
    $len = $starText.Length
    $c.Characters(1, $len).Font.Color = 255
    $ws.Range($addr).EntireRow.RowHeight = 15.75
}

# ---------------------------------------------------------------------
# 2) New block: "Satisfactory" (rows 21-22)
# ---------------------------------------------------------------------
Set-Label "A21" "Satisfactory"
Set-Value "C21" 'type: "satisfactory"'
Set-NoteRedStar "C22" "***" 'If the vaule is blank, please set  satisfactory_vUnsatisfactory=Very Unsatisfactory, satisfactory_Unsatisfactory=Unsatisfactory, satisfactory_Neutral=Neutral, satisfactory_Satisfactory=Satisfactory, satisfactory_vSatisfactory=Very Satisfactory'

# ---------------------------------------------------------------------
# 3) New block: "Image Radio" (rows 24-27)
# ---------------------------------------------------------------------
Set-Label "A24" "Image Radio"
Set-Value "C24" 'type: "img_radio"'
Set-NoteRedStar "C25" "*** " 'If the vaule is blank, please set as "NULL"'
Set-NoteRedStar "C26" "***" ' If the vaule is not blank, that means got URL, then save the value into database.'
$ws.Range("C27").EntireRow.RowHeight = 15.75
$ws.Range("C27").Font.Name = "Calibri"
$ws.Range("C27").Font.Size = 12

# ---------------------------------------------------------------------
# 4) New block: "Image Checkbox" (rows 28-32)
# ---------------------------------------------------------------------
Set-Label "A28" "Image Checkbox"
Set-Value "C28" 'type:"img_checkbox"'
Set-NoteRedStar "C29" "*** " 'If the vaule is blank, please set as "NULL"'
Set-NoteRedStar "C30" "***" ' If the vaule is not blank, that means got URL, then save the value into database.'
$ws.Range("C31").EntireRow.RowHeight = 15.75
$ws.Range("C31").Font.Name = "Calibri"
$ws.Range("C31").Font.Size = 12
$ws.Range("C32").EntireRow.RowHeight = 15.75
$ws.Range("C32").Font.Name = "Calibri"
$ws.Range("C32").Font.Size = 12

# ---------------------------------------------------------------------
# 5) Extend the "Type" -> type-string mapping table at the bottom with
#    the three new field types (rows 70, 72, 75 - mirroring the gaps
#    used by the existing entries in that table).
# ---------------------------------------------------------------------
Set-Label "A70" "Satisfactory"
Set-Value "C70" "satisfactory"

Set-Label "A72" "Image Radio"
Set-Value "C72" "img_radio"

Set-Label "A75" "Image Checkbox"
Set-Value "C75" "img_checkbox"

# ---------------------------------------------------------------------
# 6) Sheet view bookkeeping to match the saved state.
# ---------------------------------------------------------------------
$ws.Range("F77").Select() | Out-Null
